# Generate Report for Handback
# Adds "Latest Target File" / "Latest Handback File" hyperlinked values for the
# two language sheets, records the actual handback datetime, and updates the
# status text from "Ready for handoff" to "Handed back: in sync with en-US".

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------------
# 1. Status column: every cell that currently reads "Ready for handoff" now
#    reads "Handed back: in sync with en-US" (Overview!B2:C3, and Status (C)
#    column on both language sheets).
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------------
# 2. zh-cn sheet: populate "Latest Target File" (F) and "Latest Handback
#    File" (G) for both data rows, mirroring the md / xlf hyperlinks already
#    present in columns A and D.
# ---------------------------------------------------------------------------
$zhCnMdUrl2 = "https://github.com/OpenLocalizationTest/oltest/blob/cbd64e8e52aea26cc739077ff5563bc241ae9711/e2e/95350963-d258-4712-ad8f-fc017fb12334.md"
$zhCnXlfUrl2 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c99244ba16823ee6ab05b03a9f48754bab615a03/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/95350963-d258-4712-ad8f-fc017fb12334.57a8c1aae103014c88ded2631b668fa271a16ea0.zh-cn.xlf"
$zhCnMdUrl3 = "https://github.com/OpenLocalizationTest/oltest/blob/cbd64e8e52aea26cc739077ff5563bc241ae9711/e2e/b1a0afd3-03b9-4d5c-850c-92a94d0b8941.md"
$zhCnXlfUrl3 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c99244ba16823ee6ab05b03a9f48754bab615a03/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/b1a0afd3-03b9-4d5c-850c-92a94d0b8941.3299af1478829c62a486bc49036d66947efc3060.zh-cn.xlf"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("F2"), $zhCnMdUrl2, "", "", "95350963-d258-4712-ad8f-fc017fb12334.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("G2"), $zhCnXlfUrl2, "", "", "95350963-d258-4712-ad8f-fc017fb12334.57a8c1aae103014c88ded2631b668fa271a16ea0.zh-cn.xlf")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("F3"), $zhCnMdUrl3, "", "", "b1a0afd3-03b9-4d5c-850c-92a94d0b8941.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("G3"), $zhCnXlfUrl3, "", "", "b1a0afd3-03b9-4d5c-850c-92a94d0b8941.3299af1478829c62a486bc49036d66947efc3060.zh-cn.xlf")

foreach ($addr in @("F2", "G2", "F3", "G3")) {
    $wsZhCn.Range($addr).Font.Underline = 2
    $wsZhCn.Range($addr).Font.Color = 15570276
}

# Handback datetime actually recorded now (was the zero-date placeholder).
$wsZhCn.Range("H2").Value = "2016-03-17 14:38:29"
$wsZhCn.Range("H3").Value = "2016-03-17 14:38:29"

# ---------------------------------------------------------------------------
# 3. de-de sheet: same shape of change as zh-cn.
# ---------------------------------------------------------------------------
$deDeMdUrl2 = "https://github.com/OpenLocalizationTest/oltest/blob/cbd64e8e52aea26cc739077ff5563bc241ae9711/e2e/95350963-d258-4712-ad8f-fc017fb12334.md"
$deDeXlfUrl2 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f1007e9aa2bc3a9e5b92f420cf4d6f119c1fd7ee/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/95350963-d258-4712-ad8f-fc017fb12334.57a8c1aae103014c88ded2631b668fa271a16ea0.de-de.xlf"
$deDeMdUrl3 = "https://github.com/OpenLocalizationTest/oltest/blob/cbd64e8e52aea26cc739077ff5563bc241ae9711/e2e/b1a0afd3-03b9-4d5c-850c-92a94d0b8941.md"
$deDeXlfUrl3 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f1007e9aa2bc3a9e5b92f420cf4d6f119c1fd7ee/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/b1a0afd3-03b9-4d5c-850c-92a94d0b8941.3299af1478829c62a486bc49036d66947efc3060.de-de.xlf"

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("F2"), $deDeMdUrl2, "", "", "95350963-d258-4712-ad8f-fc017fb12334.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("G2"), $deDeXlfUrl2, "", "", "95350963-d258-4712-ad8f-fc017fb12334.57a8c1aae103014c88ded2631b668fa271a16ea0.de-de.xlf")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("F3"), $deDeMdUrl3, "", "", "b1a0afd3-03b9-4d5c-850c-92a94d0b8941.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("G3"), $deDeXlfUrl3, "", "", "b1a0afd3-03b9-4d5c-850c-92a94d0b8941.3299af1478829c62a486bc49036d66947efc3060.de-de.xlf")

foreach ($addr in @("F2", "G2", "F3", "G3")) {
    $wsDeDe.Range($addr).Font.Underline = 2
    $wsDeDe.Range($addr).Font.Color = 15570276
}

# Handback datetime actually recorded now (was the zero-date placeholder).
$wsDeDe.Range("H2").Value = "2016-03-17 14:38:35"
$wsDeDe.Range("H3").Value = "2016-03-17 14:38:35"
